# Updated tests in main.
# Rebuilds the "Copies served by Dep. Clerk ..." line so that:
#  - the blank-line total grows from 6 to 9 underscores, now split across
#    three runs joined by soft hyphens (mirrors the author's manual edit),
#  - the stray <w:tab/> before "Prosecutor's Office" is dropped,
#  - the "{{ defendant.first_name }}" placeholder's opening brace is fixed
#    from a single "{" to a proper "{{" and the word is spell-checked as a
#    single run instead of being split mid-word, and
#  - the "_GoBack" bookmark is moved from the final empty paragraph to sit
#    immediately before "___ Prosecutor's Office" on this line.

$d = $word.ActiveDocument

# The document already carries a hidden "_GoBack" bookmark on the trailing
# empty paragraph; remove it so it can be re-inserted at its new location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the paragraph that starts with the distinctive "Copies served..."
# text without relying on hard-coded character offsets.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Copies served by Dep. Clerk", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Copies served by Dep. Clerk' paragraph."
}

$targetPara = $searchRange.Paragraphs(1)
$targetRange = $d.Range($targetPara.Range.Start, $targetPara.Range.End - 1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Copies served by Dep. Clerk _____</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:softHyphen/></w:r><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:softHyphen/><w:t>___</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">_ on: </w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>___ Prosecutor’s Office, ___ {{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>defendant.first_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> }} {{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>defendant.last_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetRange.InsertXML($xml)
